$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps its original text representation
# (Excel would otherwise auto-convert numeric-looking strings to numbers,
# dropping formatting such as trailing zeros).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value2 = '21.652.32'
$ws.Range("E2").Value2 = '  -1.69%  '
$ws.Range("D3").Value2 = '1.532.94'
$ws.Range("E3").Value2 = '  -1.47%  '
$ws.Range("D4").Value2 = '1.001'
$ws.Range("E4").Value2 = '  +0.11%  '
$ws.Range("E5").Value2 = '  +0.14%  '
$ws.Range("D6").Value2 = '288.14'
$ws.Range("E6").Value2 = '  +0.37%  '
$ws.Range("D7").Value2 = '0.3952'
$ws.Range("E7").Value2 = '  +2.45%  '
$ws.Range("D8").Value2 = '0.3147'
$ws.Range("E8").Value2 = '  -3.07%  '
$ws.Range("D9").Value2 = '42.38'
$ws.Range("E9").Value2 = '  +1.98%  '
$ws.Range("D10").Value2 = '0.07137'
$ws.Range("E10").Value2 = '  -2.53%  '
$ws.Range("D11").Value2 = '1.041'
$ws.Range("E11").Value2 = '  -7.22%  '
$ws.Range("E12").Value2 = '  +0.11%  '
$ws.Range("D13").Value2 = '5.623'
$ws.Range("E13").Value2 = '  -1.55%  '
$ws.Range("D14").Value2 = '18.45'
$ws.Range("E14").Value2 = '  -4.60%  '
$ws.Range("D15").Value2 = '6.569'
$ws.Range("E15").Value2 = '  -3.34%  '
$ws.Range("D16").Value2 = '1.535.14'
$ws.Range("E16").Value2 = '  -1.11%  '
$ws.Range("D17").Value2 = '0.00001084'
$ws.Range("E17").Value2 = '  -0.60%  '
$ws.Range("D18").Value2 = '0.06582'
$ws.Range("E18").Value2 = '  -0.65%  '
$ws.Range("D19").Value2 = '82.95'
$ws.Range("E19").Value2 = '  -2.50%  '
$ws.Range("D20").Value2 = '1.000'
$ws.Range("E20").Value2 = '  +0.14%  '
$ws.Range("D21").Value2 = '6.088'
$ws.Range("E21").Value2 = '  -4.79%  '
$ws.Range("D22").Value2 = '15.36'
$ws.Range("E22").Value2 = '  -3.66%  '
$ws.Range("D23").Value2 = '10.80'
$ws.Range("E23").Value2 = '  -5.68%  '
$ws.Range("D24").Value2 = '2.380'
$ws.Range("E24").Value2 = '  +2.50%  '
$ws.Range("D25").Value2 = '21.649.63'
$ws.Range("E25").Value2 = '  -1.72%  '
$ws.Range("D26").Value2 = '2.336'
$ws.Range("E26").Value2 = '  -7.80%  '
$ws.Range("D27").Value2 = '147.54'
$ws.Range("E27").Value2 = '  -1.08%  '
$ws.Range("D28").Value2 = '18.25'
$ws.Range("E28").Value2 = '  -3.05%  '
$ws.Range("D29").Value2 = '4.845'
$ws.Range("E29").Value2 = '  -0.10%  '
$ws.Range("D30").Value2 = '1.709.09'
$ws.Range("E30").Value2 = '  -1.05%  '
$ws.Range("D31").Value2 = '116.59'
$ws.Range("E31").Value2 = '  -3.33%  '
$ws.Range("D32").Value2 = '5.838'
$ws.Range("E32").Value2 = '  -0.58%  '
$ws.Range("D33").Value2 = '0.9349'
$ws.Range("E33").Value2 = '  -14.46%  '
$ws.Range("D34").Value2 = '0.08134'
$ws.Range("E34").Value2 = '  -0.14%  '
$ws.Range("D35").Value2 = '8.388'
$ws.Range("E35").Value2 = '  -9.07%  '
$ws.Range("D36").Value2 = '0.06012'
$ws.Range("E36").Value2 = '  -2.90%  '
$ws.Range("E37").Value2 = '  -3.01%  '
$ws.Range("D38").Value2 = '0.02195'
$ws.Range("E38").Value2 = '  -4.36%  '
$ws.Range("D39").Value2 = '1.437'
$ws.Range("E39").Value2 = '  -13.64%  '
$ws.Range("D40").Value2 = '0.2007'
$ws.Range("E40").Value2 = '  -4.65%  '
$ws.Range("D41").Value2 = '1.172'
$ws.Range("E41").Value2 = '  -3.93%  '
$ws.Range("E42").Value2 = '  +0.12%  '
$ws.Range("D43").Value2 = '10.83'
$ws.Range("E43").Value2 = '  -0.28%  '
$ws.Range("D44").Value2 = '0.5724'
$ws.Range("E44").Value2 = '  -3.44%  '
$ws.Range("D47").Value2 = '0.5460'
$ws.Range("E47").Value2 = '  -4.84%  '
$ws.Range("D48").Value2 = '1.159'
$ws.Range("E48").Value2 = '  +0.39%  '
$ws.Range("D49").Value2 = '115.71'
$ws.Range("E49").Value2 = '  -2.98%  '
$ws.Range("D50").Value2 = '1.854'
$ws.Range("E50").Value2 = '  -3.95%  '
$ws.Range("D51").Value2 = '0.06682'
$ws.Range("E51").Value2 = '  -2.90%  '

# Rows 45 and 46 swapped positions (EnergySwap now ranks above PancakeSwap)
$ws.Range("B45").Value2 = 'EnergySwap'
$ws.Range("C45").Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value2 = '12.99'
$ws.Range("E45").Value2 = '  -3.92%  '

$ws.Range("B46").Value2 = 'PancakeSwap'
$ws.Range("C46").Value2 = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").Value2 = '3.716'
$ws.Range("E46").Value2 = '  -0.06%  '
